# Update the Sample Annotation workbook to match the new version of
# MSTemplate_Creator: the workbook should now contain only the
# "Sample_Annot" worksheet. The "Transition_Name_Annot" and "ISTD_Annot"
# worksheets (and the now-unused shared strings that only they used) are
# removed entirely.

$wb = $excel.ActiveWorkbook

# Suppress the "data will be permanently deleted" confirmation prompt that
# Excel normally raises when Delete()-ing a non-empty worksheet.
$excel.DisplayAlerts = $false

[void]$wb.Worksheets.Item("Transition_Name_Annot").Delete()
[void]$wb.Worksheets.Item("ISTD_Annot").Delete()

# Make sure the one remaining sheet is the active/selected one so the
# workbook doesn't keep a stale "active tab" pointing past the end of the
# (now much shorter) sheet list.
$wb.Worksheets.Item("Sample_Annot").Activate()

$excel.DisplayAlerts = $true
